# Remove 100.0 from conflict matrix
#
# The "100" traffic-light/phase ID is being removed from the workbook:
#   - "Setup" sheet: the row holding ID 100 (row 30) is deleted, and the
#     min-green-time value for ID 99 (row 29, col B) is changed to 0.
#   - "ConflictMatrix" sheet: both the row AND the column that correspond
#     to ID 100 (row 30 / column AE) are deleted, since the sheet is a
#     symmetric ID x ID matrix.
#   - Selection / active sheet bookkeeping is updated to match what Excel
#     leaves behind after performing these edits interactively.

$wb = $excel.ActiveWorkbook

$wsSetup = $wb.Worksheets.Item("Setup")
$wsMatrix = $wb.Worksheets.Item("ConflictMatrix")

# --- Setup sheet -----------------------------------------------------
# Update the min green time for id 99 (row 29) before the shift below.
$wsSetup.Range("B29").Value = 0

# Delete the entire row belonging to id 100 (row 30); rows below shift up.
$wsSetup.Rows.Item(30).Delete()

# --- ConflictMatrix sheet --------------------------------------------
# Delete the row for id 100 (row 30); rows below shift up.
$wsMatrix.Range("30:30").Delete()

# Delete the column for id 100 (column AE); columns to the right shift left.
$wsMatrix.Range("AE:AE").Delete()

# --- Selection / active sheet bookkeeping ----------------------------
[void]$wsSetup.Range("B30").Select()
$wsMatrix.Activate()
